$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 599, shifting the existing rows (old 599-674)
# down to become rows 601-676.
$ws.Rows.Item(599).Insert()
$ws.Rows.Item(599).Insert()

# Populate the first new row (599) with the new weekly record.
$ws.Cells.Item(599, 1).Value = 9
$ws.Cells.Item(599, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(599, 3).Value = "Metropolitana"
$ws.Cells.Item(599, 4).Value = 45154
$ws.Cells.Item(599, 5).Value = 13
$ws.Cells.Item(599, 6).Value = 100112013
$ws.Cells.Item(599, 7).Value = "Alcachofa"
$ws.Cells.Item(599, 8).Value = "Española"
$ws.Cells.Item(599, 9).Value = "Extra"
$ws.Cells.Item(599, 10).Value = 52
$ws.Cells.Item(599, 11).Value = 16000
$ws.Cells.Item(599, 12).Value = 16000
$ws.Cells.Item(599, 13).Value = 16000
$ws.Cells.Item(599, 14).Value = "$/caja 25 unidades"
$ws.Cells.Item(599, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(599, 16).Value = 16000
$ws.Cells.Item(599, 17).Value = 1
$ws.Cells.Item(599, 18).Value = "Hortaliza"

# Populate the second new row (600) with the new weekly record.
$ws.Cells.Item(600, 1).Value = 9
$ws.Cells.Item(600, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(600, 3).Value = "Metropolitana"
$ws.Cells.Item(600, 4).Value = 45154
$ws.Cells.Item(600, 5).Value = 13
$ws.Cells.Item(600, 6).Value = 100112013
$ws.Cells.Item(600, 7).Value = "Alcachofa"
$ws.Cells.Item(600, 8).Value = "Española"
$ws.Cells.Item(600, 9).Value = "Primera"
$ws.Cells.Item(600, 10).Value = 70
$ws.Cells.Item(600, 11).Value = 14000
$ws.Cells.Item(600, 12).Value = 15000
$ws.Cells.Item(600, 13).Value = 14500
$ws.Cells.Item(600, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(600, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(600, 16).Value = 483
$ws.Cells.Item(600, 17).Value = 30
$ws.Cells.Item(600, 18).Value = "Hortaliza"
